$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'41.543.36"
$ws.Range("E2").Value = "  -1.18%  "

$ws.Range("D3").Value = "'2.205.56"
$ws.Range("E3").Value = "  -0.94%  "

$ws.Range("D5").Value = "'255.35"
$ws.Range("E5").Value = "  +4.61%  "

$ws.Range("D6").Value = "'0.627"
$ws.Range("E6").Value = "  +0.07%  "

$ws.Range("D7").Value = "'69.35"
$ws.Range("E7").Value = "  +1.29%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").Value = "'0.594"
$ws.Range("E9").Value = "  +8.19%  "

$ws.Range("D10").Value = "'38.65"
$ws.Range("E10").Value = "  +8.55%  "

$ws.Range("D11").Value = "'0.0953"
$ws.Range("E11").Value = "  -0.84%  "

$ws.Range("D12").Value = "'58.37"
$ws.Range("E12").Value = "  +0.33%  "

$ws.Range("D13").Value = "'7.19"
$ws.Range("E13").Value = "  +6.97%  "

$ws.Range("D14").Value = "'0.104"
$ws.Range("E14").Value = "  -1.22%  "

$ws.Range("D15").Value = "'2.539.80"
$ws.Range("E15").Value = "  -0.71%  "

$ws.Range("D16").Value = "'14.83"
$ws.Range("E16").Value = "  +0.25%  "

$ws.Range("D17").Value = "'0.884"
$ws.Range("E17").Value = "  +4.12%  "

$ws.Range("D18").Value = "'2.199.36"
$ws.Range("E18").Value = "  -1.15%  "

$ws.Range("D19").Value = "'41.505.09"
$ws.Range("E19").Value = "  -1.07%  "

$ws.Range("D20").Value = "'0.0₃0956"
$ws.Range("E20").Value = "  +0.04%  "

$ws.Range("D21").Value = "'6.25"
$ws.Range("E21").Value = "  +1.94%  "

$ws.Range("D22").Value = "'72.21"
$ws.Range("E22").Value = "  -0.47%  "

$ws.Range("D23").Value = "'233.21"
$ws.Range("E23").Value = "  -0.66%  "

$ws.Range("E24").Value = "  +1.41%  "

$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").Value = "'11.88"
$ws.Range("E25").Value = "  +20.20%  "

$ws.Range("B26").Value = "WEMIXToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D26").Value = "'3.88"
$ws.Range("E26").Value = "  +7.03%  "

$ws.Range("E27").Value = "  +0.01%  "

$ws.Range("D28").Value = "'2.56"
$ws.Range("E28").Value = "  +4.65%  "

$ws.Range("E29").Value = "  -2.06%  "

$ws.Range("D30").Value = "'170.98"
$ws.Range("E30").Value = "  -0.23%  "

$ws.Range("D31").Value = "'20.67"
$ws.Range("E31").Value = "  +1.16%  "

$ws.Range("D32").Value = "'0.122"
$ws.Range("E32").Value = "  +1.00%  "

$ws.Range("D33").Value = "'5.55"
$ws.Range("E33").Value = "  +6.71%  "

$ws.Range("D34").Value = "'0.123"
$ws.Range("E34").Value = "  -1.53%  "

$ws.Range("D35").Value = "'0.0735"
$ws.Range("E35").Value = "  +2.97%  "

$ws.Range("D36").Value = "'26.22"
$ws.Range("E36").Value = "  +16.47%  "

$ws.Range("D37").Value = "'4.64"
$ws.Range("E37").Value = "  -0.32%  "

$ws.Range("D38").Value = "'3.99"
$ws.Range("E38").Value = "  +2.30%  "

$ws.Range("D39").Value = "'0.0301"
$ws.Range("E39").Value = "  +8.27%  "

$ws.Range("E40").Value = "  -2.34%  "

$ws.Range("D41").Value = "'5.84"
$ws.Range("E41").Value = "  +0.20%  "

$ws.Range("D42").Value = "'12.03"
$ws.Range("E42").Value = "  +19.22%  "

$ws.Range("D43").Value = "'64.19"
$ws.Range("E43").Value = "  -3.59%  "

$ws.Range("B44").Value = "FTXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D44").Value = "'4.98"
$ws.Range("E44").Value = "  +0.28%  "

$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").Value = "'0.204"
$ws.Range("E45").Value = "  +6.45%  "

$ws.Range("D46").Value = "'8.68"
$ws.Range("E46").Value = "  -3.62%  "

$ws.Range("D47").Value = "'0.101"
$ws.Range("E47").Value = "  +0.32%  "

$ws.Range("D48").Value = "'1.01"
$ws.Range("E48").Value = "  +0.42%  "

$ws.Range("D49").Value = "'1.16"
$ws.Range("E49").Value = "  +4.36%  "

$ws.Range("D50").Value = "'4.39"
$ws.Range("E50").Value = "  -3.73%  "

$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").Value = "'1.18"
$ws.Range("E51").Value = "  -0.29%  "
